$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Add new row 7 with the new changelog entry
$ws.Range("A7").Value = "2014-29-01"
$ws.Range("C7").Value = "não"
$ws.Range("B7").Value = "Removida as colunas de valor e data para cadastro das despesas/receitas. Adicionado todos os meses e ano"

# Widen column B to fit the new longer text (98 characters)
$ws.Columns.Item(2).ColumnWidth = 97.16666666666667

# Update the active selection as recorded in the sheet view
$ws.Range("B11").Select()
